# Remove first-person language from the incident report: the sentence
# "... a network protocol analyzer, I discovered that there were a very
# large number of TCP SYN requests coming from an unfamiliar IP address."
# becomes "... a network protocol analyzer, the cybersecurity analyst
# discovered that there were a very large number of TCP SYN requests
# coming from an unfamiliar IP address."

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "analyzer, I discovered",   # FindText
    $true,                      # MatchCase
    $false,                     # MatchWholeWord
    $false,                     # MatchWildcards
    $false,                     # MatchSoundsLike
    $false,                     # MatchAllWordForms
    $true,                      # Forward
    1,                          # Wrap (wdFindContinue)
    $false,                     # Format
    "analyzer, the cybersecurity analyst discovered",  # ReplaceWith
    2                           # Replace (wdReplaceAll)
)
